# Insert a new data row at row 410 (shifts existing rows 410-520 down to 411-521)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 410; this shifts rows 410:520 -> 411:521
$ws.Rows("410:410").Insert()

# Populate the newly inserted row 410 with the new record.
$ws.Cells.Item(410, 1).Value = 11
$ws.Cells.Item(410, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(410, 3).Value = "Bíobío"
$ws.Cells.Item(410, 4).Value = 44785
$ws.Cells.Item(410, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(410, 5).Value = 8
$ws.Cells.Item(410, 6).Value = "Fruta"
$ws.Cells.Item(410, 7).Value = 100104
$ws.Cells.Item(410, 8).Value = "Frutos de pepita"
$ws.Cells.Item(410, 9).Value = 100104005
$ws.Cells.Item(410, 10).Value = "Pera"
$ws.Cells.Item(410, 11).Value = "Packham's Triumph"
$ws.Cells.Item(410, 12).Value = "Primera"
$ws.Cells.Item(410, 13).Value = 220
$ws.Cells.Item(410, 14).Value = 8500
$ws.Cells.Item(410, 15).Value = 9000
$ws.Cells.Item(410, 16).Value = 8773
$ws.Cells.Item(410, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(410, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(410, 19).Value = 548
$ws.Cells.Item(410, 20).Value = 16
